$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price values remain stored as text (matching original formatting)
$textCells = @("D5", "D6", "D9", "D12", "D13", "D14", "D17", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D42", "D46", "D49", "D50", "D51")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "65.137.96"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("D3").Value = "3.180.06"
$ws.Range("E3").Value = "  +3.88%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "573.78"
$ws.Range("E5").Value = "  +2.56%  "
$ws.Range("D6").Value = "151.19"
$ws.Range("E6").Value = "  +5.01%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.175.83"
$ws.Range("E8").Value = "  +3.83%  "
$ws.Range("D9").Value = "0.528"
$ws.Range("E9").Value = "  +3.13%  "
$ws.Range("E10").Value = "  +4.46%  "
$ws.Range("E11").Value = "  +2.08%  "
$ws.Range("D12").Value = "0.506"
$ws.Range("E12").Value = "  +4.99%  "
$ws.Range("D13").Value = "0.0000277"
$ws.Range("E13").Value = "  +19.40%  "
$ws.Range("D14").Value = "38.22"
$ws.Range("E14").Value = "  +7.00%  "
$ws.Range("D15").Value = "3.691.55"
$ws.Range("E15").Value = "  +3.60%  "
$ws.Range("D16").Value = "65.221.30"
$ws.Range("E16").Value = "  +1.68%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").Value = "7.22"
$ws.Range("E17").Value = "  +6.56%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.171.78"
$ws.Range("E18").Value = "  +3.39%  "
$ws.Range("E19").Value = "  +1.19%  "
$ws.Range("D20").Value = "512.64"
$ws.Range("E20").Value = "  +7.13%  "
$ws.Range("D21").Value = "14.96"
$ws.Range("E21").Value = "  +6.76%  "
$ws.Range("D22").Value = "0.738"
$ws.Range("E22").Value = "  +8.32%  "
$ws.Range("D23").Value = "15.72"
$ws.Range("E23").Value = "  +9.33%  "
$ws.Range("D24").Value = "7.85"
$ws.Range("E24").Value = "  +3.08%  "
$ws.Range("D25").Value = "84.99"
$ws.Range("E25").Value = "  +3.43%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "9.22"
$ws.Range("E27").Value = "  +15.08%  "
$ws.Range("D28").Value = "2.91"
$ws.Range("E28").Value = "  +3.79%  "
$ws.Range("D29").Value = "2.21"
$ws.Range("E29").Value = "  +8.03%  "
$ws.Range("D30").Value = "28.16"
$ws.Range("E30").Value = "  +6.78%  "
$ws.Range("D31").Value = "2.81"
$ws.Range("E31").Value = "  +15.15%  "
$ws.Range("D32").Value = "1.23"
$ws.Range("E32").Value = "  +7.55%  "
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("D34").Value = "6.35"
$ws.Range("E34").Value = "  +11.46%  "
$ws.Range("D35").Value = "6.69"
$ws.Range("E35").Value = "  +7.07%  "
$ws.Range("D36").Value = "55.66"
$ws.Range("E36").Value = "  +1.26%  "
$ws.Range("D37").Value = "482.26"
$ws.Range("E37").Value = "  +8.04%  "
$ws.Range("D38").Value = "0.0892"
$ws.Range("E38").Value = "  +9.66%  "
$ws.Range("D39").Value = "3.12"
$ws.Range("E39").Value = "  +8.67%  "
$ws.Range("D40").Value = "0.0423"
$ws.Range("E40").Value = "  +3.24%  "
$ws.Range("D41").Value = "3.136.97"
$ws.Range("E41").Value = "  +4.82%  "
$ws.Range("D42").Value = "8.65"
$ws.Range("E42").Value = "  +4.74%  "
$ws.Range("E43").Value = "  +4.17%  "
$ws.Range("E44").Value = "  +17.76%  "
$ws.Range("E45").Value = "  +10.76%  "
$ws.Range("D46").Value = "29.36"
$ws.Range("E46").Value = "  +4.67%  "
$ws.Range("D47").Value = "0.0₃0596"
$ws.Range("E47").Value = "  +14.52%  "
$ws.Range("D49").Value = "0.116"
$ws.Range("E49").Value = "  +2.06%  "
$ws.Range("D50").Value = "2.31"
$ws.Range("E50").Value = "  +11.28%  "
$ws.Range("D51").Value = "122.62"
$ws.Range("E51").Value = "  +2.92%  "

# Reset style back to Normal for cells where we forced text format, to avoid leftover numeric style on the cell
foreach ($ref in $textCells) {
    $ws.Range($ref).Style = "Normal"
}
